$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "SREDNIA"
$ws.Range("D12").Formula = "=AVERAGE(D2:D11)"
$ws.Range("D12").NumberFormat = "General"

$ws.Range("C13").Value = "G"
$ws.Range("D13").Formula = "=(4*3.141*3.141*0.396)/(D12*D12)"

$ws.Range("C14").Value = "BLAD"
$ws.Range("D14").Formula = "=(D13-9.81)/9.81 * 100"

$ws.Range("C15").Select()
